$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every Coin/Link/Price/Volume cell in this sheet is stored as literal
# text (even numeric-looking prices like "13.40"), never as a real
# number, so each write below is apostrophe-prefixed. That is the COM
# equivalent of a user typing an '-prefixed entry into the cell: Excel
# stores the text verbatim (no silent numeric coercion / dropped
# trailing zeros) while the leading apostrophe itself is not part of
# the stored value.

$ws.Range("D2").Value = "'42.554.88"
$ws.Range("E2").Value = "'  -0.42%  "

$ws.Range("D3").Value = "'2.290.85"
$ws.Range("E3").Value = "'  -0.36%  "

$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "'  +0.65%  "

$ws.Range("D5").Value = "'311.73"
$ws.Range("E5").Value = "'  -3.45%  "

$ws.Range("D6").Value = "'102.71"
$ws.Range("E6").Value = "'  -1.92%  "

$ws.Range("D7").Value = "'0.623"
$ws.Range("E7").Value = "'  -1.18%  "

$ws.Range("E8").Value = "'  -0.11%  "

$ws.Range("E9").Value = "'  -1.07%  "

$ws.Range("D10").Value = "'38.88"
$ws.Range("E10").Value = "'  -3.33%  "

$ws.Range("E11").Value = "'  -1.36%  "

$ws.Range("E12").Value = "'  -3.88%  "

$ws.Range("D14").Value = "'0.978"
$ws.Range("E14").Value = "'  +0.32%  "

$ws.Range("D15").Value = "'15.27"
$ws.Range("E15").Value = "'  -0.29%  "

$ws.Range("D16").Value = "'2.637.11"
$ws.Range("E16").Value = "'  -0.38%  "

$ws.Range("D17").Value = "'2.287.42"
$ws.Range("E17").Value = "'  -0.02%  "

$ws.Range("D18").Value = "'42.720.09"
$ws.Range("E18").Value = "'  +0.17%  "

$ws.Range("D19").Value = "'7.28"
$ws.Range("E19").Value = "'  -3.08%  "

$ws.Range("E20").Value = "'  -2.20%  "

$ws.Range("D21").Value = "'13.40"
$ws.Range("E21").Value = "'  -0.45%  "

$ws.Range("D22").Value = "'73.34"
$ws.Range("E22").Value = "'  -0.39%  "

$ws.Range("D23").Value = "'268.68"
$ws.Range("E23").Value = "'  -0.83%  "

$ws.Range("E24").Value = "'  -6.29%  "

$ws.Range("D25").Value = "'2.20"
$ws.Range("E25").Value = "'  -1.77%  "

$ws.Range("D26").Value = "'1.01"
$ws.Range("E26").Value = "'  +0.22%  "

$ws.Range("B27").Value = "'Filecoin"
$ws.Range("C27").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D27").Value = "'7.26"
$ws.Range("E27").Value = "'  +17.40%  "

$ws.Range("B28").Value = "'Cosmos"
$ws.Range("C28").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'10.75"
$ws.Range("E28").Value = "'  -1.58%  "

$ws.Range("D29").Value = "'2.29"
$ws.Range("E29").Value = "'  -1.43%  "

$ws.Range("E30").Value = "'  -1.11%  "

$ws.Range("D31").Value = "'35.64"
$ws.Range("E31").Value = "'  -7.00%  "

$ws.Range("D32").Value = "'164.24"

$ws.Range("E33").Value = "'  -4.01%  "

$ws.Range("D34").Value = "'0.129"
$ws.Range("E34").Value = "'  -2.08%  "

$ws.Range("E35").Value = "'  +1.10%  "

$ws.Range("E36").Value = "'  -2.88%  "

$ws.Range("D37").Value = "'4.52"
$ws.Range("E37").Value = "'  -2.29%  "

$ws.Range("D38").Value = "'0.0346"
$ws.Range("E38").Value = "'  -2.70%  "

$ws.Range("D39").Value = "'2.79"
$ws.Range("E39").Value = "'  +1.99%  "

$ws.Range("D40").Value = "'3.60"
$ws.Range("E40").Value = "'  -3.66%  "

$ws.Range("D41").Value = "'107.19"
$ws.Range("E41").Value = "'  +8.25%  "

$ws.Range("D42").Value = "'1.56"
$ws.Range("E42").Value = "'  +0.62%  "

$ws.Range("D43").Value = "'70.93"
$ws.Range("E43").Value = "'  +0.84%  "

$ws.Range("B44").Value = "'FirstDigitalUSD"
$ws.Range("C44").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.01"
$ws.Range("E44").Value = "'  +0.25%  "

$ws.Range("B45").Value = "'Algorand"
$ws.Range("C45").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "'0.226"
$ws.Range("E45").Value = "'  +0.48%  "

$ws.Range("D46").Value = "'12.04"
$ws.Range("E46").Value = "'  -3.04%  "

$ws.Range("D47").Value = "'1.730.98"
$ws.Range("E47").Value = "'  +8.54%  "

$ws.Range("D48").Value = "'110.41"
$ws.Range("E48").Value = "'  -2.63%  "

$ws.Range("D49").Value = "'76.36"
$ws.Range("E49").Value = "'  -7.38%  "

$ws.Range("E50").Value = "'  -2.79%  "

$ws.Range("D51").Value = "'5.14"
$ws.Range("E51").Value = "'  -2.59%  "
